# 0607: optimize constraint (3)
# The parameter-sweep table on Sheet1 shrinks from 10 scenario rows to 4,
# and the surviving rows get new constraint values (H1 tolerance 0.05/0.01/0.00
# -> 0.02; the bench-weight cap column collapses to just 0/80; a new "E"
# exposure-bound column takes 0.5/1.0; and mkt_type now alternates
# CSI500/CSI300). The backing Excel Table ("表1") auto-resizes with the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the bottom 5 scenario rows (rows 6-10) first; the ListObject backing
# the range shrinks (and re-anchors to A1:P5) automatically when the rows
# underneath it are removed.
$ws.Range("A6:P10").EntireRow.Delete()

# Re-write the 4 remaining data rows (2-5) with the new scenario values.
# Columns: A run | B alpha_name | C mkt_type | D beta_kind | E beta_suffix |
# F beta_args | G H0 | H H1 | I B | J E | K D | L N | M wei_tole |
# N begin_date | O end_date | P opt_verbose

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "FRtn5D(0.0,3.0)"
$ws.Range("C2").Value = "CSI500"
$ws.Range("D2").Value = "Barra"
$ws.Range("E2").Value = "barra3"
$ws.Range("F2").Value = "(['size', 'beta', 'momentum'],)"
$ws.Range("G2").Value = "0.20"
$ws.Range("H2").Value = "0.02"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0.5"
$ws.Range("K2").Value = "2"
$ws.Range("L2").Value = "inf"
$ws.Range("M2").Value = "1e-5"
$ws.Range("N2").Value = "2016-02-01"
$ws.Range("O2").Value = "2022-03-31"
$ws.Range("P2").Value = "'FALSE"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "FRtn5D(0.0,3.0)"
$ws.Range("C3").Value = "CSI300"
$ws.Range("D3").Value = "Barra"
$ws.Range("E3").Value = "barra3"
$ws.Range("F3").Value = "(['size', 'beta', 'momentum'],)"
$ws.Range("G3").Value = "0.20"
$ws.Range("H3").Value = "0.02"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "1.0"
$ws.Range("K3").Value = "2"
$ws.Range("L3").Value = "inf"
$ws.Range("M3").Value = "1e-5"
$ws.Range("N3").Value = "2016-02-01"
$ws.Range("O3").Value = "2022-03-31"
$ws.Range("P3").Value = "'FALSE"

$ws.Range("A4").Value = "1"
$ws.Range("B4").Value = "FRtn5D(0.0,3.0)"
$ws.Range("C4").Value = "CSI500"
$ws.Range("D4").Value = "Barra"
$ws.Range("E4").Value = "barra3"
$ws.Range("F4").Value = "(['size', 'beta', 'momentum'],)"
$ws.Range("G4").Value = "0.20"
$ws.Range("H4").Value = "0.02"
$ws.Range("I4").Value = "80"
$ws.Range("J4").Value = "0.5"
$ws.Range("K4").Value = "2"
$ws.Range("L4").Value = "inf"
$ws.Range("M4").Value = "1e-5"
$ws.Range("N4").Value = "2016-02-01"
$ws.Range("O4").Value = "2022-03-31"
$ws.Range("P4").Value = "'FALSE"

$ws.Range("A5").Value = "1"
$ws.Range("B5").Value = "FRtn5D(0.0,3.0)"
$ws.Range("C5").Value = "CSI300"
$ws.Range("D5").Value = "Barra"
$ws.Range("E5").Value = "barra3"
$ws.Range("F5").Value = "(['size', 'beta', 'momentum'],)"
$ws.Range("G5").Value = "0.20"
$ws.Range("H5").Value = "0.02"
$ws.Range("I5").Value = "80"
$ws.Range("J5").Value = "1.0"
$ws.Range("K5").Value = "2"
$ws.Range("L5").Value = "inf"
$ws.Range("M5").Value = "1e-5"
$ws.Range("N5").Value = "2016-02-01"
$ws.Range("O5").Value = "2022-03-31"
$ws.Range("P5").Value = "'FALSE"

# Match the author's post-edit selection/cursor position.
$ws.Range("B7").Select()
